# Update CDA Logical model for ST.r2b
#
# The "Metadata" sheet gains a new "Jurisdiction" property row (with an
# empty value) right after the "Contact" row, and the "Version" /
# "Date" values are bumped. The "Elements" sheet content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row for "Jurisdiction" right after the "Contact" row (row 10),
# copying the formatting of the existing data rows so the new row matches
# the sheet's established style.
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Bump the published Version and Date metadata values.
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
